# Adds a new "2017" inflow/otflow column pair to the district migration
# table, shifting 2018-2022 data two columns to the right (C:D -> E:F, etc.)
# and appending a new 2022 column pair at M:N.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 ("inflow"/"otflow" labels): extend the alternating pattern through M2:N2
$ws.Range("C2").HorizontalAlignment = -4108
$ws.Range("C2").Font.Bold = $true
$ws.Range("C2").Value = "inflow"
$ws.Range("D2").HorizontalAlignment = -4108
$ws.Range("D2").Font.Bold = $true
$ws.Range("D2").Value = "otflow"
$ws.Range("E2").HorizontalAlignment = -4108
$ws.Range("E2").Font.Bold = $true
$ws.Range("E2").Value = "inflow"
$ws.Range("F2").HorizontalAlignment = -4108
$ws.Range("F2").Font.Bold = $true
$ws.Range("F2").Value = "otflow"
$ws.Range("G2").HorizontalAlignment = -4108
$ws.Range("G2").Font.Bold = $true
$ws.Range("G2").Value = "inflow"
$ws.Range("H2").HorizontalAlignment = -4108
$ws.Range("H2").Font.Bold = $true
$ws.Range("H2").Value = "otflow"
$ws.Range("I2").HorizontalAlignment = -4108
$ws.Range("I2").Font.Bold = $true
$ws.Range("I2").Value = "inflow"
$ws.Range("J2").HorizontalAlignment = -4108
$ws.Range("J2").Font.Bold = $true
$ws.Range("J2").Value = "otflow"
$ws.Range("K2").HorizontalAlignment = -4108
$ws.Range("K2").Font.Bold = $true
$ws.Range("K2").Value = "inflow"
$ws.Range("L2").HorizontalAlignment = -4108
$ws.Range("L2").Font.Bold = $true
$ws.Range("L2").Value = "otflow"
$ws.Range("M2").HorizontalAlignment = -4108
$ws.Range("M2").Font.Bold = $true
$ws.Range("M2").Value = "inflow"
$ws.Range("N2").HorizontalAlignment = -4108
$ws.Range("N2").Font.Bold = $true
$ws.Range("N2").Value = "otflow"

# Row 3 (year labels): 2017 now occupies C3:D3, existing years shift right to E3:N3
$ws.Range("C3").HorizontalAlignment = -4108
$ws.Range("C3").Font.Bold = $true
$ws.Range("C3").Value = 2017
$ws.Range("D3").HorizontalAlignment = -4108
$ws.Range("D3").Font.Bold = $true
$ws.Range("D3").Value = 2017
$ws.Range("E3").HorizontalAlignment = -4108
$ws.Range("E3").Font.Bold = $true
$ws.Range("E3").Value = 2018
$ws.Range("F3").HorizontalAlignment = -4108
$ws.Range("F3").Font.Bold = $true
$ws.Range("F3").Value = 2018
$ws.Range("G3").HorizontalAlignment = -4108
$ws.Range("G3").Font.Bold = $true
$ws.Range("G3").Value = 2019
$ws.Range("H3").HorizontalAlignment = -4108
$ws.Range("H3").Font.Bold = $true
$ws.Range("H3").Value = 2019
$ws.Range("I3").HorizontalAlignment = -4108
$ws.Range("I3").Font.Bold = $true
$ws.Range("I3").Value = 2020
$ws.Range("J3").HorizontalAlignment = -4108
$ws.Range("J3").Font.Bold = $true
$ws.Range("J3").Value = 2020
$ws.Range("K3").HorizontalAlignment = -4108
$ws.Range("K3").Font.Bold = $true
$ws.Range("K3").Value = 2021
$ws.Range("L3").HorizontalAlignment = -4108
$ws.Range("L3").Font.Bold = $true
$ws.Range("L3").Value = 2021
$ws.Range("M3").HorizontalAlignment = -4108
$ws.Range("M3").Font.Bold = $true
$ws.Range("M3").Value = 2022
$ws.Range("N3").HorizontalAlignment = -4108
$ws.Range("N3").Font.Bold = $true
$ws.Range("N3").Value = 2022

# Row 4 (total LO): 2017 data in C4:D4, existing years shift right to E4:N4
$ws.Range("C4").HorizontalAlignment = -4108
$ws.Range("C4").Value = 105704
$ws.Range("D4").HorizontalAlignment = -4108
$ws.Range("D4").Value = 74845
$ws.Range("E4").HorizontalAlignment = -4108
$ws.Range("E4").Value = 124177
$ws.Range("F4").HorizontalAlignment = -4108
$ws.Range("F4").Value = 80436
$ws.Range("G4").HorizontalAlignment = -4108
$ws.Range("G4").Value = 103076
$ws.Range("H4").HorizontalAlignment = -4108
$ws.Range("H4").Value = 68599
$ws.Range("I4").HorizontalAlignment = -4108
$ws.Range("I4").Value = 90215
$ws.Range("J4").HorizontalAlignment = -4108
$ws.Range("J4").Value = 62860
$ws.Range("K4").HorizontalAlignment = -4108
$ws.Range("K4").Value = 95323
$ws.Range("L4").HorizontalAlignment = -4108
$ws.Range("L4").Value = 65345
$ws.Range("M4").HorizontalAlignment = -4108
$ws.Range("M4").Value = 94031
$ws.Range("N4").HorizontalAlignment = -4108
$ws.Range("N4").Value = 64066

# Row 5 (Central dist.): 2017 data in C5:D5, existing years shift right to E5:N5
$ws.Range("C5").HorizontalAlignment = -4108
$ws.Range("C5").Value = 7821
$ws.Range("D5").HorizontalAlignment = -4108
$ws.Range("D5").Value = 6172
$ws.Range("E5").HorizontalAlignment = -4108
$ws.Range("E5").Value = 9608
$ws.Range("F5").HorizontalAlignment = -4108
$ws.Range("F5").Value = 5978
$ws.Range("G5").HorizontalAlignment = -4108
$ws.Range("G5").Value = 8413
$ws.Range("H5").HorizontalAlignment = -4108
$ws.Range("H5").Value = 5352
$ws.Range("I5").HorizontalAlignment = -4108
$ws.Range("I5").Value = 7302
$ws.Range("J5").HorizontalAlignment = -4108
$ws.Range("J5").Value = 4995
$ws.Range("K5").HorizontalAlignment = -4108
$ws.Range("K5").Value = 7094
$ws.Range("L5").HorizontalAlignment = -4108
$ws.Range("L5").Value = 4806
$ws.Range("M5").HorizontalAlignment = -4108
$ws.Range("M5").Value = 6929
$ws.Range("N5").HorizontalAlignment = -4108
$ws.Range("N5").Value = 5047

# Row 6 (North-west dist.): 2017 data in C6:D6, existing years shift right to E6:N6
$ws.Range("C6").HorizontalAlignment = -4108
$ws.Range("C6").Value = 58153
$ws.Range("D6").HorizontalAlignment = -4108
$ws.Range("D6").Value = 49440
$ws.Range("E6").HorizontalAlignment = -4108
$ws.Range("E6").Value = 69556
$ws.Range("F6").HorizontalAlignment = -4108
$ws.Range("F6").Value = 51953
$ws.Range("G6").HorizontalAlignment = -4108
$ws.Range("G6").Value = 64646
$ws.Range("H6").HorizontalAlignment = -4108
$ws.Range("H6").Value = 50775
$ws.Range("I6").HorizontalAlignment = -4108
$ws.Range("I6").Value = 57076
$ws.Range("J6").HorizontalAlignment = -4108
$ws.Range("J6").Value = 45421
$ws.Range("K6").HorizontalAlignment = -4108
$ws.Range("K6").Value = 63328
$ws.Range("L6").HorizontalAlignment = -4108
$ws.Range("L6").Value = 48569
$ws.Range("M6").HorizontalAlignment = -4108
$ws.Range("M6").Value = 62611
$ws.Range("N6").HorizontalAlignment = -4108
$ws.Range("N6").Value = 47291

# Row 7 (South dist.): 2017 data in C7:D7, existing years shift right to E7:N7
$ws.Range("C7").HorizontalAlignment = -4108
$ws.Range("C7").Value = 4018
$ws.Range("D7").HorizontalAlignment = -4108
$ws.Range("D7").Value = 2634
$ws.Range("E7").HorizontalAlignment = -4108
$ws.Range("E7").Value = 5143
$ws.Range("F7").HorizontalAlignment = -4108
$ws.Range("F7").Value = 2795
$ws.Range("G7").HorizontalAlignment = -4108
$ws.Range("G7").Value = 4991
$ws.Range("H7").HorizontalAlignment = -4108
$ws.Range("H7").Value = 2748
$ws.Range("I7").HorizontalAlignment = -4108
$ws.Range("I7").Value = 4325
$ws.Range("J7").HorizontalAlignment = -4108
$ws.Range("J7").Value = 2870
$ws.Range("K7").HorizontalAlignment = -4108
$ws.Range("K7").Value = 4174
$ws.Range("L7").HorizontalAlignment = -4108
$ws.Range("L7").Value = 2784
$ws.Range("M7").HorizontalAlignment = -4108
$ws.Range("M7").Value = 4362
$ws.Range("N7").HorizontalAlignment = -4108
$ws.Range("N7").Value = 2649

# Row 8 (North-Caucas dist.): 2017 data in C8:D8, existing years shift right to E8:N8
$ws.Range("C8").HorizontalAlignment = -4108
$ws.Range("C8").Value = 2713
$ws.Range("D8").HorizontalAlignment = -4108
$ws.Range("D8").Value = 1300
$ws.Range("E8").HorizontalAlignment = -4108
$ws.Range("E8").Value = 3182
$ws.Range("F8").HorizontalAlignment = -4108
$ws.Range("F8").Value = 1473
$ws.Range("G8").HorizontalAlignment = -4108
$ws.Range("G8").Value = 3065
$ws.Range("H8").HorizontalAlignment = -4108
$ws.Range("H8").Value = 1453
$ws.Range("I8").HorizontalAlignment = -4108
$ws.Range("I8").Value = 2649
$ws.Range("J8").HorizontalAlignment = -4108
$ws.Range("J8").Value = 1343
$ws.Range("K8").HorizontalAlignment = -4108
$ws.Range("K8").Value = 2328
$ws.Range("L8").HorizontalAlignment = -4108
$ws.Range("L8").Value = 1291
$ws.Range("M8").HorizontalAlignment = -4108
$ws.Range("M8").Value = 2313
$ws.Range("N8").HorizontalAlignment = -4108
$ws.Range("N8").Value = 1432

# Row 9 (Volga dist.): 2017 data in C9:D9, existing years shift right to E9:N9
$ws.Range("C9").HorizontalAlignment = -4108
$ws.Range("C9").Value = 7659
$ws.Range("D9").HorizontalAlignment = -4108
$ws.Range("D9").Value = 3586
$ws.Range("E9").HorizontalAlignment = -4108
$ws.Range("E9").Value = 9618
$ws.Range("F9").HorizontalAlignment = -4108
$ws.Range("F9").Value = 3706
$ws.Range("G9").HorizontalAlignment = -4108
$ws.Range("G9").Value = 9107
$ws.Range("H9").HorizontalAlignment = -4108
$ws.Range("H9").Value = 3784
$ws.Range("I9").HorizontalAlignment = -4108
$ws.Range("I9").Value = 7546
$ws.Range("J9").HorizontalAlignment = -4108
$ws.Range("J9").Value = 3686
$ws.Range("K9").HorizontalAlignment = -4108
$ws.Range("K9").Value = 7081
$ws.Range("L9").HorizontalAlignment = -4108
$ws.Range("L9").Value = 3582
$ws.Range("M9").HorizontalAlignment = -4108
$ws.Range("M9").Value = 7017
$ws.Range("N9").HorizontalAlignment = -4108
$ws.Range("N9").Value = 3359

# Row 10 (Ural dist.): 2017 data in C10:D10, existing years shift right to E10:N10
$ws.Range("C10").HorizontalAlignment = -4108
$ws.Range("C10").Value = 2713
$ws.Range("D10").HorizontalAlignment = -4108
$ws.Range("D10").Value = 1233
$ws.Range("E10").HorizontalAlignment = -4108
$ws.Range("E10").Value = 3669
$ws.Range("F10").HorizontalAlignment = -4108
$ws.Range("F10").Value = 1291
$ws.Range("G10").HorizontalAlignment = -4108
$ws.Range("G10").Value = 3663
$ws.Range("H10").HorizontalAlignment = -4108
$ws.Range("H10").Value = 1301
$ws.Range("I10").HorizontalAlignment = -4108
$ws.Range("I10").Value = 3122
$ws.Range("J10").HorizontalAlignment = -4108
$ws.Range("J10").Value = 1376
$ws.Range("K10").HorizontalAlignment = -4108
$ws.Range("K10").Value = 3038
$ws.Range("L10").HorizontalAlignment = -4108
$ws.Range("L10").Value = 1321
$ws.Range("M10").HorizontalAlignment = -4108
$ws.Range("M10").Value = 2932
$ws.Range("N10").HorizontalAlignment = -4108
$ws.Range("N10").Value = 1313

# Row 11 (Siberia dist.): 2017 data in C11:D11, existing years shift right to E11:N11
$ws.Range("C11").HorizontalAlignment = -4108
$ws.Range("C11").Value = 4973
$ws.Range("D11").HorizontalAlignment = -4108
$ws.Range("D11").Value = 1811
$ws.Range("E11").HorizontalAlignment = -4108
$ws.Range("E11").Value = 5303
$ws.Range("F11").HorizontalAlignment = -4108
$ws.Range("F11").Value = 1827
$ws.Range("G11").HorizontalAlignment = -4108
$ws.Range("G11").Value = 5272
$ws.Range("H11").HorizontalAlignment = -4108
$ws.Range("H11").Value = 1816
$ws.Range("I11").HorizontalAlignment = -4108
$ws.Range("I11").Value = 4724
$ws.Range("J11").HorizontalAlignment = -4108
$ws.Range("J11").Value = 1864
$ws.Range("K11").HorizontalAlignment = -4108
$ws.Range("K11").Value = 4650
$ws.Range("L11").HorizontalAlignment = -4108
$ws.Range("L11").Value = 1689
$ws.Range("M11").HorizontalAlignment = -4108
$ws.Range("M11").Value = 4440
$ws.Range("N11").HorizontalAlignment = -4108
$ws.Range("N11").Value = 1703

# Row 12 (Far-east dist.): 2017 data in C12:D12, existing years shift right to E12:N12
$ws.Range("C12").HorizontalAlignment = -4108
$ws.Range("C12").Value = 2625
$ws.Range("D12").HorizontalAlignment = -4108
$ws.Range("D12").Value = 1036
$ws.Range("E12").HorizontalAlignment = -4108
$ws.Range("E12").Value = 4175
$ws.Range("F12").HorizontalAlignment = -4108
$ws.Range("F12").Value = 1258
$ws.Range("G12").HorizontalAlignment = -4108
$ws.Range("G12").Value = 3919
$ws.Range("H12").HorizontalAlignment = -4108
$ws.Range("H12").Value = 1370
$ws.Range("I12").HorizontalAlignment = -4108
$ws.Range("I12").Value = 3471
$ws.Range("J12").HorizontalAlignment = -4108
$ws.Range("J12").Value = 1305
$ws.Range("K12").HorizontalAlignment = -4108
$ws.Range("K12").Value = 3630
$ws.Range("L12").HorizontalAlignment = -4108
$ws.Range("L12").Value = 1303
$ws.Range("M12").HorizontalAlignment = -4108
$ws.Range("M12").Value = 3427
$ws.Range("N12").HorizontalAlignment = -4108
$ws.Range("N12").Value = 1272

# Restore the active-cell selection recorded in the workbook
$ws.Range("D14").Select() | Out-Null
